# "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
# The account-statement periods (Periodo Mora, column E, rows 16-20) are
# refreshed: the previous batch of periods (1803-1807, oldest-to-newest) is
# replaced by re-entering the same periods in newest-to-oldest order, which
# is how the underlying macro regenerates this sheet from the updated
# database each time it is run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E16").Value = "1807"
$ws.Range("E17").Value = "1806"
$ws.Range("E18").Value = "1805"
$ws.Range("E19").Value = "1804"
$ws.Range("E20").Value = "1803"
